# Insert a new row at position 124 (shifts existing rows 124-192 down to 125-193)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly data point
$ws.Range("A124").Value = 7
$ws.Range("B124").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C124").Value = 'Ñuble'
$ws.Range("D124").Value = 44518
$ws.Range("E124").Value = 16
$ws.Range("F124").Value = 100112023
$ws.Range("G124").Value = 'Brócoli'
$ws.Range("H124").Value = 'Sin especificar'
$ws.Range("I124").Value = 'Primera'
$ws.Range("J124").Value = 300
$ws.Range("K124").Value = 650
$ws.Range("L124").Value = 700
$ws.Range("M124").Value = 675
$ws.Range("N124").Value = '$/unidad'
$ws.Range("O124").Value = 'Región del Maule'
$ws.Range("P124").Value = 675
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = 'Hortaliza'
